$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.457.05"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.519.64"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'539.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'139.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "2.524.70"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "'5.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "'0.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "2.966.32"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'23.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "59.287.48"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "2.515.76"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'11.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "'326.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'0.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'7.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "'6.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").Value = "0.0₃0782"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").Value = "'165.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").Value = "'18.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "'36.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("D42").Value = "'5.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").Value = "'280.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.77%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'0.601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'10.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "'123.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'17.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.82%  "
